$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price (D) and Volume(1h) (E) columns to Text format so numeric-looking
# strings (e.g. "1.001", "0.9993", "0.000006810") are not auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "26.170.74"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").Value = "1.748.26"
$ws.Range("E3").Value = "  +0.20%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "242.36"
$ws.Range("E5").Value = "  +4.01%  "
$ws.Range("D6").Value = "0.5537"
$ws.Range("E6").Value = "  +6.78%  "
$ws.Range("D7").Value = "0.9993"
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("D8").Value = "0.2835"
$ws.Range("E8").Value = "  +0.77%  "
$ws.Range("D9").Value = "0.06182"
$ws.Range("E9").Value = "  +1.04%  "
$ws.Range("D10").Value = "1.757.24"
$ws.Range("E10").Value = "  +0.68%  "
$ws.Range("D11").Value = "0.07189"
$ws.Range("E11").Value = "  +1.67%  "
$ws.Range("D12").Value = "15.56"
$ws.Range("E12").Value = "  +1.36%  "
$ws.Range("D13").Value = "0.6642"
$ws.Range("E13").Value = "  +3.44%  "
$ws.Range("D14").Value = "4.658"
$ws.Range("E14").Value = "  +2.95%  "
$ws.Range("D15").Value = "78.39"
$ws.Range("E15").Value = "  +1.33%  "
$ws.Range("D16").Value = "0.9997"
$ws.Range("E16").Value = "  -0.08%  "
$ws.Range("D17").Value = "1.001"
$ws.Range("E17").Value = "  -0.08%  "
$ws.Range("D18").Value = "26.117.64"
$ws.Range("E18").Value = "  +0.54%  "
$ws.Range("D19").Value = "11.86"
$ws.Range("E19").Value = "  +3.16%  "
$ws.Range("D20").Value = "0.000006810"
$ws.Range("E20").Value = "  +3.11%  "
$ws.Range("D21").Value = "1.986.67"
$ws.Range("E21").Value = "  +0.39%  "
$ws.Range("D22").Value = "4.406"
$ws.Range("E22").Value = "  +6.16%  "
$ws.Range("D23").Value = "8.807"
$ws.Range("E23").Value = "  +1.84%  "
$ws.Range("D24").Value = "5.316"
$ws.Range("E24").Value = "  +3.34%  "
$ws.Range("D25").Value = "140.17"
$ws.Range("E25").Value = "  +0.40%  "
$ws.Range("D26").Value = "1.522"
$ws.Range("E26").Value = "  +0.72%  "
$ws.Range("D28").Value = "1.818"
$ws.Range("E28").Value = "  -0.16%  "
$ws.Range("D29").Value = "105.80"
$ws.Range("E29").Value = "  +2.94%  "
$ws.Range("D30").Value = "0.08527"
$ws.Range("E30").Value = "  +3.09%  "
$ws.Range("D31").Value = "3.796"
$ws.Range("E31").Value = "  +3.51%  "
$ws.Range("D32").Value = "3.647"
$ws.Range("E32").Value = "  +5.90%  "
$ws.Range("D33").Value = "0.04648"
$ws.Range("E33").Value = "  +3.78%  "
$ws.Range("D34").Value = "2.654"
$ws.Range("E34").Value = "  +1.49%  "
$ws.Range("D36").Value = "0.6294"
$ws.Range("E36").Value = "  +1.85%  "
$ws.Range("D37").Value = "2.710"
$ws.Range("E37").Value = "  +1.26%  "
$ws.Range("D38").Value = "0.01616"
$ws.Range("E38").Value = "  +1.68%  "
$ws.Range("D39").Value = "1.979"
$ws.Range("E39").Value = "  +2.97%  "
$ws.Range("D40").Value = "0.9985"
$ws.Range("E40").Value = "  -0.21%  "
$ws.Range("D41").Value = "99.78"
$ws.Range("E41").Value = "  -0.40%  "
$ws.Range("D42").Value = "0.3947"
$ws.Range("E42").Value = "  +2.41%  "
$ws.Range("D43").Value = "0.7526"
$ws.Range("E43").Value = "  +1.87%  "
$ws.Range("D44").Value = "5.041"
$ws.Range("E44").Value = "  -0.52%  "
$ws.Range("D45").Value = "0.1153"
$ws.Range("E45").Value = "  +2.11%  "
$ws.Range("D46").Value = "6.363"
$ws.Range("E46").Value = "  +0.49%  "
$ws.Range("D47").Value = "0.05344"
$ws.Range("E47").Value = "  -2.06%  "
$ws.Range("D48").Value = "55.03"
$ws.Range("E48").Value = "  +3.63%  "
$ws.Range("D49").Value = "31.00"
$ws.Range("E49").Value = "  +2.89%  "
$ws.Range("D50").Value = "0.3511"
$ws.Range("E50").Value = "  +2.67%  "
$ws.Range("D51").Value = "7.679"
$ws.Range("E51").Value = "  +0.26%  "

# Rows with only Volume(1h) change, Price stays the same
$ws.Range("E27").Value = "  +1.61%  "
$ws.Range("E35").Value = "  +1.99%  "

# Restore default cell style (no explicit style index) to match original formatting
$ws.Range("D2:E51").Style = "Normal"
